$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.61"
$ws.Range("E2").Value = "'-0.47%"
$ws.Range("D3").Value = "'37.08"
$ws.Range("E3").Value = "'5.88%"
$ws.Range("D4").Value = "'5.005"
$ws.Range("E4").Value = "'-3.27%"
$ws.Range("E5").Value = "'0.54%"
$ws.Range("D6").Value = "'2.229"
$ws.Range("E6").Value = "'-2.88%"
$ws.Range("D7").Value = "'7.998"
$ws.Range("E7").Value = "'-0.76%"
$ws.Range("D8").Value = "'4.015"
$ws.Range("E8").Value = "'0.72%"
$ws.Range("D9").Value = "'0.9199"
$ws.Range("E9").Value = "'-0.39%"
$ws.Range("D10").Value = "'0.09545"
$ws.Range("E10").Value = "'-4.96%"
$ws.Range("D11").Value = "'0.1884"
$ws.Range("E11").Value = "'2.30%"
$ws.Range("D12").Value = "'0.08551"
$ws.Range("E12").Value = "'0.47%"
$ws.Range("D13").Value = "'0.03586"
$ws.Range("E13").Value = "'6.39%"
$ws.Range("D14").Value = "'0.09973"
$ws.Range("E14").Value = "'0.56%"
$ws.Range("D15").Value = "'0.001482"
$ws.Range("E15").Value = "'-0.22%"
$ws.Range("D16").Value = "'0.005713"
$ws.Range("E16").Value = "'0.45%"
$ws.Range("E17").Value = "'-0.74%"
$ws.Range("D18").Value = "'2.255"
$ws.Range("E18").Value = "'5.95%"
$ws.Range("E19").Value = "'-0.63%"
$ws.Range("D20").Value = "'0.1316"
$ws.Range("E20").Value = "'-0.77%"
$ws.Range("D21").Value = "'4.759"
$ws.Range("E21").Value = "'3.82%"
$ws.Range("E22").Value = "'-8.16%"
$ws.Range("D23").Value = "'0.04601"
$ws.Range("E23").Value = "'-1.08%"
$ws.Range("D24").Value = "'0.001231"
$ws.Range("E24").Value = "'0.70%"
$ws.Range("D25").Value = "'0.004460"
$ws.Range("E25").Value = "'3.10%"
$ws.Range("E26").Value = "'7.53%"
$ws.Range("E27").Value = "'39.71%"
$ws.Range("D39").Value = "'0.01810"
$ws.Range("D40").Value = "'0.04725"
$ws.Range("E40").Value = "'-0.50%"
$ws.Range("D41").Value = "'0.008112"
$ws.Range("E41").Value = "'5.35%"
$ws.Range("D42").Value = "'0.1396"
$ws.Range("E42").Value = "'-1.27%"
$ws.Range("D43").Value = "'0.007545"
$ws.Range("E43").Value = "'6.51%"
$ws.Range("D44").Value = "'0.002228"
$ws.Range("E44").Value = "'-2.77%"
$ws.Range("D45").Value = "'0.01044"
$ws.Range("E45").Value = "'3.79%"
$ws.Range("D46").Value = "'0.00006169"
$ws.Range("E46").Value = "'2.57%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.10%"
$ws.Range("D48").Value = "'0.0005795"
$ws.Range("E48").Value = "'-0.10%"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("D51").Value = "'0.00002098"
$ws.Range("E51").Value = "'-0.10%"
